$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '57.088.67'
$ws.Range('E2').Value = '  +4.59%  '
Set-TextValue $ws.Range('D3') '3.247.41'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '396.02'
$ws.Range('E5').Value = '  -1.35%  '
Set-TextValue $ws.Range('D6') '108.13'
$ws.Range('E6').Value = '  -1.62%  '
Set-TextValue $ws.Range('D7') '0.589'
$ws.Range('E7').Value = '  +7.12%  '
Set-TextValue $ws.Range('D8') '3.244.09'
$ws.Range('E8').Value = '  +2.46%  '
$ws.Range('E9').Value = '  +0.05%  '
Set-TextValue $ws.Range('D11') '39.25'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('E12').Value = '  +9.92%  '
$ws.Range('E13').Value = '  +2.06%  '
Set-TextValue $ws.Range('D14') '3.759.09'
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('E16').Value = '  +0.05%  '
Set-TextValue $ws.Range('D17') '3.251.62'
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('E18').Value = '  -2.48%  '
$ws.Range('E19').Value = '  +1.79%  '
Set-TextValue $ws.Range('D20') '56.963.03'
Set-TextValue $ws.Range('D21') '3.33'
$ws.Range('E21').Value = '  +0.99%  '
$ws.Range('E22').Value = '  +7.83%  '
$ws.Range('E23').Value = '  +1.11%  '
Set-TextValue $ws.Range('D24') '294.06'
$ws.Range('E24').Value = '  +6.80%  '
Set-TextValue $ws.Range('D25') '74.16'
Set-TextValue $ws.Range('D26') '3.16'
$ws.Range('E26').Value = '  -2.94%  '
Set-TextValue $ws.Range('D27') '28.08'
$ws.Range('E27').Value = '  +0.99%  '
Set-TextValue $ws.Range('D29') '7.64'
$ws.Range('E29').Value = '  -5.29%  '
Set-TextValue $ws.Range('D30') '7.20'
$ws.Range('E30').Value = '  -6.05%  '
$ws.Range('E31').Value = '  -1.44%  '
Set-TextValue $ws.Range('D33') '11.21'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('E35').Value = '  +8.64%  '
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('E38').Value = '  +0.61%  '
Set-TextValue $ws.Range('D39') '0.999'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -4.49%  '
Set-TextValue $ws.Range('D41') '2.96'
$ws.Range('E41').Value = '  +2.92%  '
Set-TextValue $ws.Range('D42') '139.84'
$ws.Range('E42').Value = '  +6.08%  '
$ws.Range('E43').Value = '  +3.98%  '
$ws.Range('E44').Value = '  -1.67%  '
Set-TextValue $ws.Range('D45') '17.07'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('E46').Value = '  -3.79%  '
Set-TextValue $ws.Range('D47') '0.279'
$ws.Range('E47').Value = '  -4.49%  '
$ws.Range('E48').Value = '  +11.52%  '
Set-TextValue $ws.Range('D49') '22.17'
$ws.Range('E49').Value = '  +0.32%  '
Set-TextValue $ws.Range('D50') '2.160.58'
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('E51').Value = '  -5.69%  '
